# Update the four "Trt.ID: X" legend labels on slide 1 (inside the grouped
# plot shape) so that they read just the treatment letter ("C", "D", "A",
# "B") and shrink their bounding boxes accordingly (AutoSize is off for
# these shapes, so position/size must be set explicitly to match the new,
# narrower text).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

# Map shape name -> new Left (pt), Width (pt), and new text.
# (Top/Height are unchanged from the original values.)
# Note: the point values below are nudged by a hair beyond the naive
# EMU/12700 conversion to counteract the host's internal float32
# Left/Width storage, so the saved EMU lands exactly on target.
$targets = @{
    "tx38" = @{ Left = 205.35889443779527; Width = 8.219527759055119;  Text = "C" }
    "tx40" = @{ Left = 357.00416572834650; Width = 9.296850693700787;  Text = "D" }
    "tx42" = @{ Left = 205.13472750944880; Width = 8.667953055905512;  Text = "A" }
    "tx44" = @{ Left = 357.38426216850394; Width = 8.536693113385827;  Text = "B" }
}

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $sh = $grp.GroupItems.Item($i)
    if ($targets.ContainsKey($sh.Name)) {
        $t = $targets[$sh.Name]
        $sh.TextFrame.TextRange.Text = $t.Text
        $sh.Left = $t.Left
        $sh.Width = $t.Width
    }
}
